$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-5 from 45204 to 45207
$ws.Range("C2:C5").Value = 45207
